$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.372558333333333
$ws.Range("H2").Value = 4.117675
$ws.Range("I2").Value = 0.3340102211301095
$ws.Range("J2").Value = 0.3340102211301095
$ws.Range("M2").Value = 14.25737566666667
$ws.Range("N2").Value = 42.772127
$ws.Range("O2").Value = 0.2087950866344732
$ws.Range("P2").Value = 0.2087950866344732
$ws.Range("Q2").Value = 19.56907978274722
$ws.Range("R2").Value = 176.121718044725
$ws.Range("S2").Value = 0.06973969305766077
$ws.Range("T2").Value = 0.06973969305766077
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.372558333333333
$ws.Range("H3").Value = 4.117675
$ws.Range("I3").Value = 0.3340102211301095
$ws.Range("J3").Value = 0.3340102211301095
$ws.Range("N3").Value = 87.128332
$ws.Range("O3").Value = 0.4253229592313036
$ws.Range("P3").Value = 0.4253229592313036
$ws.Range("Q3").Value = 39.86290605201111
$ws.Range("R3").Value = 358.7661544681
$ws.Range("S3").Value = 0.1420622156645603
$ws.Range("T3").Value = 0.1420622156645603
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.372558333333333
$ws.Range("H4").Value = 4.117675
$ws.Range("I4").Value = 0.3340102211301095
$ws.Range("J4").Value = 0.3340102211301095
$ws.Range("M4").Value = 20.11084633333333
$ws.Range("N4").Value = 60.332539
$ws.Range("O4").Value = 0.2945174484164121
$ws.Range("P4").Value = 0.2945174484164122
$ws.Range("Q4").Value = 27.60330972520277
$ws.Range("R4").Value = 248.429787526825
$ws.Range("S4").Value = 0.09837183807224145
$ws.Range("T4").Value = 0.09837183807224147
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.372558333333333
$ws.Range("H5").Value = 4.117675
$ws.Range("I5").Value = 0.3340102211301095
$ws.Range("J5").Value = 0.3340102211301095
$ws.Range("M5").Value = 4.873057999999999
$ws.Range("N5").Value = 14.619174
$ws.Range("O5").Value = 0.07136450571781097
$ws.Range("P5").Value = 0.07136450571781099
$ws.Range("Q5").Value = 6.688556366716666
$ws.Range("R5").Value = 60.19700730045
$ws.Range("S5").Value = 0.02383647433564701
$ws.Range("T5").Value = 0.02383647433564701
$ws.Range("I6").Value = 0.01293592767872722
$ws.Range("J6").Value = 0.01293592767872721
$ws.Range("M6").Value = 14.25737566666667
$ws.Range("N6").Value = 42.772127
$ws.Range("O6").Value = 0.2087950866344732
$ws.Range("P6").Value = 0.2087950866344732
$ws.Range("Q6").Value = 0.7578935756886668
$ws.Range("R6").Value = 6.821042181198001
$ws.Range("S6").Value = 0.002700958140377129
$ws.Range("T6").Value = 0.002700958140377129
$ws.Range("I7").Value = 0.01293592767872722
$ws.Range("J7").Value = 0.01293592767872721
$ws.Range("N7").Value = 87.128332
$ws.Range("O7").Value = 0.4253229592313036
$ws.Range("P7").Value = 0.4253229592313036
$ws.Range("S7").Value = 0.005501947040718388
$ws.Range("T7").Value = 0.005501947040718388
$ws.Range("I8").Value = 0.01293592767872722
$ws.Range("J8").Value = 0.01293592767872721
$ws.Range("M8").Value = 20.11084633333333
$ws.Range("N8").Value = 60.332539
$ws.Range("O8").Value = 0.2945174484164121
$ws.Range("P8").Value = 0.2945174484164122
$ws.Range("Q8").Value = 1.069052369387333
$ws.Range("R8").Value = 9.621471324486
$ws.Range("S8").Value = 0.003809856412837981
$ws.Range("T8").Value = 0.003809856412837981
$ws.Range("I9").Value = 0.01293592767872722
$ws.Range("J9").Value = 0.01293592767872721
$ws.Range("M9").Value = 4.873057999999999
$ws.Range("N9").Value = 14.619174
$ws.Range("O9").Value = 0.07136450571781097
$ws.Range("P9").Value = 0.07136450571781099
$ws.Range("Q9").Value = 0.259042017164
$ws.Range("R9").Value = 2.331378154476
$ws.Range("S9").Value = 0.0009231660847937176
$ws.Range("T9").Value = 0.0009231660847937177
$ws.Range("G10").Value = 2.683614
$ws.Range("H10").Value = 8.050841999999999
$ws.Range("I10").Value = 0.6530538511911632
$ws.Range("J10").Value = 0.6530538511911632
$ws.Range("M10").Value = 14.25737566666667
$ws.Range("N10").Value = 42.772127
$ws.Range("O10").Value = 0.2087950866344732
$ws.Range("P10").Value = 0.2087950866344732
$ws.Range("Q10").Value = 38.261292942326
$ws.Range("R10").Value = 344.351636480934
$ws.Range("S10").Value = 0.1363544354364353
$ws.Range("T10").Value = 0.1363544354364353
$ws.Range("G11").Value = 2.683614
$ws.Range("H11").Value = 8.050841999999999
$ws.Range("I11").Value = 0.6530538511911632
$ws.Range("J11").Value = 0.6530538511911632
$ws.Range("N11").Value = 87.128332
$ws.Range("O11").Value = 0.4253229592313036
$ws.Range("P11").Value = 0.4253229592313036
$ws.Range("Q11").Value = 77.939603850616
$ws.Range("R11").Value = 701.4564346555439
$ws.Range("S11").Value = 0.2777587965260249
$ws.Range("T11").Value = 0.277758796526025
$ws.Range("G12").Value = 2.683614
$ws.Range("H12").Value = 8.050841999999999
$ws.Range("I12").Value = 0.6530538511911632
$ws.Range("J12").Value = 0.6530538511911632
$ws.Range("M12").Value = 20.11084633333333
$ws.Range("N12").Value = 60.332539
$ws.Range("O12").Value = 0.2945174484164121
$ws.Range("P12").Value = 0.2945174484164122
$ws.Range("Q12").Value = 53.96974877198199
$ws.Range("R12").Value = 485.727738947838
$ws.Range("S12").Value = 0.1923357539313327
$ws.Range("T12").Value = 0.1923357539313327
$ws.Range("G13").Value = 2.683614
$ws.Range("H13").Value = 8.050841999999999
$ws.Range("I13").Value = 0.6530538511911632
$ws.Range("J13").Value = 0.6530538511911632
$ws.Range("M13").Value = 4.873057999999999
$ws.Range("N13").Value = 14.619174
$ws.Range("O13").Value = 0.07136450571781097
$ws.Range("P13").Value = 0.07136450571781099
$ws.Range("Q13").Value = 13.077406671612
$ws.Range("R13").Value = 117.696660044508
$ws.Range("S13").Value = 0.04660486529737024
$ws.Range("T13").Value = 0.04660486529737025
